$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Total" -> "Total_Pop"
$ws.Range("B1").Value = "Total_Pop"

# Strip " County, New York" suffix from the Geographic column entries
$ws.Range("A2").Value = "Albany"
$ws.Range("A3").Value = "Allegany"
$ws.Range("A4").Value = "Bronx"
$ws.Range("A5").Value = "Broome"
$ws.Range("A6").Value = "Cattaraugus"
$ws.Range("A7").Value = "Cayuga"
$ws.Range("A8").Value = "Chautauqua"
$ws.Range("A9").Value = "Chemung"
$ws.Range("A10").Value = "Chenango"
$ws.Range("A11").Value = "Clinton"
$ws.Range("A12").Value = "Columbia"
$ws.Range("A13").Value = "Cortland"
$ws.Range("A14").Value = "Delaware"
$ws.Range("A15").Value = "Dutchess"
$ws.Range("A16").Value = "Erie"
$ws.Range("A17").Value = "Essex"
$ws.Range("A18").Value = "Franklin"
$ws.Range("A19").Value = "Fulton"
$ws.Range("A20").Value = "Genesee"
$ws.Range("A21").Value = "Greene"
$ws.Range("A22").Value = "Hamilton"
$ws.Range("A23").Value = "Herkimer"
$ws.Range("A24").Value = "Jefferson"
$ws.Range("A25").Value = "Kings"
$ws.Range("A26").Value = "Lewis"
$ws.Range("A27").Value = "Livingston"
$ws.Range("A28").Value = "Madison"
$ws.Range("A29").Value = "Monroe"
$ws.Range("A30").Value = "Montgomery"
$ws.Range("A31").Value = "Nassau"
$ws.Range("A32").Value = "New York"
$ws.Range("A33").Value = "Niagara"
$ws.Range("A34").Value = "Oneida"
$ws.Range("A35").Value = "Onondaga"
$ws.Range("A36").Value = "Ontario"
$ws.Range("A37").Value = "Orange"
$ws.Range("A38").Value = "Orleans"
$ws.Range("A39").Value = "Oswego"
$ws.Range("A40").Value = "Otsego"
$ws.Range("A41").Value = "Putnam"
$ws.Range("A42").Value = "Queens"
$ws.Range("A43").Value = "Rensselaer"
$ws.Range("A44").Value = "Richmond"
$ws.Range("A45").Value = "Rockland"
$ws.Range("A46").Value = "St. Lawrence"
$ws.Range("A47").Value = "Saratoga"
$ws.Range("A48").Value = "Schenectady"
$ws.Range("A49").Value = "Schoharie"
$ws.Range("A50").Value = "Schuyler"
$ws.Range("A51").Value = "Seneca"
$ws.Range("A52").Value = "Steuben"
$ws.Range("A53").Value = "Suffolk"
$ws.Range("A54").Value = "Sullivan"
$ws.Range("A55").Value = "Tioga"
$ws.Range("A56").Value = "Tompkins"
$ws.Range("A57").Value = "Ulster"
$ws.Range("A58").Value = "Warren"
$ws.Range("A59").Value = "Washington"
$ws.Range("A60").Value = "Wayne"
$ws.Range("A61").Value = "Westchester"
$ws.Range("A62").Value = "Wyoming"
$ws.Range("A63").Value = "Yates"

# Resize columns (matches narrower layout after the rename/merge)
$ws.Columns.Item(1).ColumnWidth = 42.333333333333336
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 8.666666666666666
$ws.Columns.Item(9).ColumnWidth = 21.5

# Select whole column A (as the author did before saving)
[void]$ws.Columns.Item(1).Select()
